$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = $origStyle
}

# Row 2
Set-TextValue $ws.Range("D2") '43.969.59'
Set-TextValue $ws.Range("E2") '  -0.03%  '

# Row 3
Set-TextValue $ws.Range("D3") '2.248.27'
Set-TextValue $ws.Range("E3") '  -1.92%  '

# Row 4
Set-TextValue $ws.Range("E4") '  -0.02%  '

# Row 5
Set-TextValue $ws.Range("D5") '231.93'
Set-TextValue $ws.Range("E5") '  -0.07%  '

# Row 6
Set-TextValue $ws.Range("D6") '0.631'
Set-TextValue $ws.Range("E6") '  +0.26%  '

# Row 7
Set-TextValue $ws.Range("D7") '63.05'
Set-TextValue $ws.Range("E7") '  -1.26%  '

# Row 8
Set-TextValue $ws.Range("E8") '  -0.10%  '

# Row 9
Set-TextValue $ws.Range("E9") '  +5.84%  '

# Row 10
Set-TextValue $ws.Range("E10") '  +3.13%  '

# Row 11
Set-TextValue $ws.Range("D11") '57.33'
Set-TextValue $ws.Range("E11") '  -0.55%  '

# Row 12
Set-TextValue $ws.Range("D12") '26.24'
Set-TextValue $ws.Range("E12") '  -0.83%  '

# Row 13
Set-TextValue $ws.Range("E13") '  +0.66%  '

# Row 14
Set-TextValue $ws.Range("D14") '2.582.25'
Set-TextValue $ws.Range("E14") '  -1.92%  '

# Row 15
Set-TextValue $ws.Range("D15") '15.50'
Set-TextValue $ws.Range("E15") '  -2.36%  '

# Row 16
Set-TextValue $ws.Range("E16") '  +2.49%  '

# Row 17
Set-TextValue $ws.Range("D17") '0.824'
Set-TextValue $ws.Range("E17") '  +0.83%  '

# Row 18
Set-TextValue $ws.Range("D18") '2.261.60'
Set-TextValue $ws.Range("E18") '  -1.01%  '

# Row 19
Set-TextValue $ws.Range("D19") '43.869.52'
Set-TextValue $ws.Range("E19") '  +0.07%  '

# Row 20
Set-TextValue $ws.Range("D20") '0.0₃0981'
Set-TextValue $ws.Range("E20") '  +2.69%  '

# Row 21
Set-TextValue $ws.Range("D21") '72.63'
Set-TextValue $ws.Range("E21") '  -0.93%  '

# Row 22
Set-TextValue $ws.Range("D22") '6.05'
Set-TextValue $ws.Range("E22") '  -2.36%  '

# Row 23
Set-TextValue $ws.Range("D23") '245.12'
Set-TextValue $ws.Range("E23") '  -2.80%  '

# Row 24
Set-TextValue $ws.Range("E24") '  -0.09%  '

# Row 25
Set-TextValue $ws.Range("E25") '  -7.53%  '

# Row 26
Set-TextValue $ws.Range("B26") 'WEMIXToken'
Set-TextValue $ws.Range("C26") 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
Set-TextValue $ws.Range("D26") '3.35'
Set-TextValue $ws.Range("E26") '  +21.46%  '

# Row 27
Set-TextValue $ws.Range("B27") 'Toncoin'
Set-TextValue $ws.Range("C27") 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
Set-TextValue $ws.Range("D27") '2.21'
Set-TextValue $ws.Range("E27") '  -11.00%  '

# Row 28
Set-TextValue $ws.Range("D28") '9.80'
Set-TextValue $ws.Range("E28") '  -1.12%  '

# Row 29
Set-TextValue $ws.Range("D29") '171.94'
Set-TextValue $ws.Range("E29") '  +0.13%  '

# Row 30
Set-TextValue $ws.Range("D30") '21.17'
Set-TextValue $ws.Range("E30") '  +2.91%  '

# Row 31
Set-TextValue $ws.Range("D31") '0.138'
Set-TextValue $ws.Range("E31") '  -1.38%  '

# Row 32
Set-TextValue $ws.Range("D32") '1.40'
Set-TextValue $ws.Range("E32") '  -3.04%  '

# Row 33
Set-TextValue $ws.Range("E33") '  +1.94%  '

# Row 34
Set-TextValue $ws.Range("E34") '  -2.30%  '

# Row 35
Set-TextValue $ws.Range("E35") '  +0.86%  '

# Row 36
Set-TextValue $ws.Range("E36") '  -4.16%  '

# Row 37
Set-TextValue $ws.Range("D37") '3.64'
Set-TextValue $ws.Range("E37") '  -1.92%  '

# Row 38
Set-TextValue $ws.Range("D38") '6.40'
Set-TextValue $ws.Range("E38") '  -3.34%  '

# Row 39
Set-TextValue $ws.Range("D39") '2.27'
Set-TextValue $ws.Range("E39") '  -4.31%  '

# Row 40
Set-TextValue $ws.Range("D40") '0.0251'
Set-TextValue $ws.Range("E40") '  +0.89%  '

# Row 41
Set-TextValue $ws.Range("E41") '  -0.06%  '

# Row 42
Set-TextValue $ws.Range("D42") '8.56'
Set-TextValue $ws.Range("E42") '  -0.21%  '

# Row 43
Set-TextValue $ws.Range("E43") '  +1.62%  '

# Row 44
Set-TextValue $ws.Range("D44") '17.02'

# Row 45
Set-TextValue $ws.Range("D45") '97.24'
Set-TextValue $ws.Range("E45") '  -1.27%  '

# Row 46
Set-TextValue $ws.Range("E46") '  -2.41%  '

# Row 47
Set-TextValue $ws.Range("E47") '  -2.69%  '

# Row 48
Set-TextValue $ws.Range("D48") '4.32'
Set-TextValue $ws.Range("E48") '  -7.63%  '

# Row 49
Set-TextValue $ws.Range("D49") '1.440.31'
Set-TextValue $ws.Range("E49") '  -3.37%  '

# Row 50
Set-TextValue $ws.Range("E50") '  -3.25%  '

# Row 51
Set-TextValue $ws.Range("E51") '  +1.80%  '
